$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns to match the latest crypto data snapshot.
# Force text number-format on the Price column so values like "1.00" or "579.94" are
# stored as literal text (matching the source feed), not auto-converted to numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.515.13"
$ws.Range("E2").Value = "  +0.78%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.446.38"
$ws.Range("E3").Value = "  +1.69%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.94"
$ws.Range("E5").Value = "  +1.03%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.83"
$ws.Range("E6").Value = "  +8.97%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.448.34"
$ws.Range("E7").Value = "  +1.94%  "

$ws.Range("E8").Value = "  +0.09%  "

$ws.Range("E9").Value = "  +0.81%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.81"
$ws.Range("E10").Value = "  +2.45%  "

$ws.Range("E11").Value = "  +1.27%  "

$ws.Range("E12").Value = "  +1.12%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.035.26"
$ws.Range("E13").Value = "  +1.88%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.91"
$ws.Range("E14").Value = "  +6.55%  "

$ws.Range("E16").Value = "  +0.39%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.446.06"
$ws.Range("E17").Value = "  +1.81%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.649.63"
$ws.Range("E18").Value = "  +0.84%  "

$ws.Range("E19").Value = "  +8.48%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.33"
$ws.Range("E20").Value = "  +2.02%  "

$ws.Range("E21").Value = "  +0.52%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "388.97"
$ws.Range("E22").Value = "  +3.12%  "

$ws.Range("E23").Value = "  +2.59%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.587.34"
$ws.Range("E24").Value = "  +1.76%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.02"
$ws.Range("E25").Value = "  +2.51%  "

$ws.Range("E26").Value = "  +0.61%  "

$ws.Range("E27").Value = "  +0.14%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000124"
$ws.Range("E28").Value = "  -1.26%  "

$ws.Range("E29").Value = "  +5.84%  "

$ws.Range("E30").Value = "  +3.47%  "

$ws.Range("E31").Value = "  +0.11%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.54"
$ws.Range("E32").Value = "  -13.11%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.26"
$ws.Range("E33").Value = "  +1.08%  "

$ws.Range("E34").Value = "  +1.02%  "

$ws.Range("E35").Value = "  -0.04%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "24.02"
$ws.Range("E36").Value = "  +1.20%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.27"
$ws.Range("E37").Value = "  +0.70%  "

$ws.Range("E38").Value = "  +2.68%  "

$ws.Range("E39").Value = "  +0.97%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "166.61"
$ws.Range("E40").Value = "  +1.63%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0792"
$ws.Range("E41").Value = "  +3.77%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.57"
$ws.Range("E42").Value = "  +9.98%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.794"
$ws.Range("E43").Value = "  +2.23%  "

$ws.Range("E44").Value = "  +1.77%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +0.10%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "42.29"
$ws.Range("E46").Value = "  +1.54%  "

$ws.Range("E47").Value = "  -0.37%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.608.75"
$ws.Range("E48").Value = "  +5.69%  "

$ws.Range("E49").Value = "  -3.12%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.05"
$ws.Range("E50").Value = "  +3.31%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.07"
$ws.Range("E51").Value = "  -0.76%  "
